$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): drop the old "Links" label in D1 ---
$ws.Range("D1").ClearContents()

# --- Column headers (row 2) ---
$ws.Range("C2").Value = "Pricing"
$ws.Range("D2").Value = "Website / Source"
$ws.Range("E2").ClearContents()

# --- Data rows: column D now shows the plain URL instead of the page title,
#     and column E ("Notes & Links") is cleared out entirely ---

# Row 3 - Agent Laboratory (D3 already held the plain URL)
$ws.Range("D3").Value = "https://agentlaboratory.github.io/"
$ws.Range("E3").ClearContents()

# Row 4 - PaSa
$ws.Range("D4").Value = "https://github.com/bytedance/pasa"
$ws.Range("E4").ClearContents()

# Row 5 - SciSciGPT
$ws.Range("D5").Value = "https://sciscigpt.com/"
$ws.Range("E5").ClearContents()

# Row 6 - Clarivate Academic Agents
$ws.Range("D6").Value = "https://clarivate.com/life-sciences-healthcare/lp/discover-drg-fusion/?campaignname=LS_DRG_Fusion_DRG_Reports_Store_LeadGen_Paid&campaignid=701VO00000Y1Y47YAF&utm_campaign=LS_DRG_Fusion_DRG_Reports_Store_LeadGen_Paid&utm_source=Google&utm_medium=Paid_Search&_bt=750691571036&_bk=clarivate%20drg%20fusion&_bm=b&_bn=g&_bg=178464204465&gad_source=1&gad_campaignid=22523827876&gclid=CjwKCAjwyb3DBhBlEiwAqZLe5PFTyFQYJocByawtmw7j1BsJsjnBB0sbYBEVC6PlP19pOPTzkapBABoC24cQAvD_BwE"
$ws.Range("E6").ClearContents()

# Row 7 - STORM
$ws.Range("D7").Value = "https://storm.genie.stanford.edu/"
$ws.Range("E7").ClearContents()

# Row 8 - NotebookLM
$ws.Range("D8").Value = "https://admin.google.com/ServiceNotAllowed?application=692380834322&source=scrip&continue=https://notebooklm.google.com/?original_referer%3Dhttps://www.google.com%2523%26pli%3D1"
$ws.Range("E8").ClearContents()

# Row 9 - Jenni AI
$ws.Range("D9").Value = "https://jenni.ai/?utm_source=google&utm_medium=cpc&utm_campaign=19905599675&utm_term=jenni%20ai&utm_content=652760871019&cmc_adid=ga_652760871019_19905599675&utm_group=146838708919&gad_source=1&gad_campaignid=19905599675&gclid=CjwKCAjwyb3DBhBlEiwAqZLe5H7gdC-v6NZBFAUpr542ZZpHgKLRpK5rR8FJ2jedJuevFdmhFMDmGxoCwmEQAvD_BwE"
$ws.Range("E9").ClearContents()

# Row 10 - Elicit
$ws.Range("D10").Value = "https://elicit.com/"
$ws.Range("E10").ClearContents()

# Row 11 - Iris.ai
$ws.Range("D11").Value = "https://iris.ai/"
$ws.Range("E11").ClearContents()

# Row 12 - Semantic Scholar
$ws.Range("D12").Value = "https://www.semanticscholar.org/"
$ws.Range("E12").ClearContents()

# Row 13 - ChatPDF
$ws.Range("D13").Value = "https://www.chatpdf.com/"
$ws.Range("E13").ClearContents()

# Row 14 - Scite.AI
$ws.Range("D14").Value = "https://scite.ai/"
$ws.Range("E14").ClearContents()

# Row 15 - ClickUp Research Agent
$ws.Range("D15").Value = "https://clickup.com/brain"
$ws.Range("E15").ClearContents()

# Row 16 - Afforai/Logically.app
$ws.Range("D16").Value = "https://afforai.com/"
$ws.Range("E16").ClearContents()

# Row 17 - Julius AI
$ws.Range("D17").Value = "https://julius.ai/"
$ws.Range("E17").ClearContents()

# Row 18 - Scifocus
$ws.Range("D18").Value = "https://www.scifocus.ai/"
$ws.Range("E18").ClearContents()

# --- View state: scroll back to top-left and move the active selection to D23 ---
$ws.Activate()
$aw = $excel.ActiveWindow
$aw.ScrollRow = 1
$aw.ScrollColumn = 1
$ws.Range("D23").Select()
